$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before C (old C,D,E -> new E,F,G), making room for "modelo" and "politica"
$ws.Columns("C:D").Insert()

# Header row
$ws.Range("C1").Value = "modelo"
$ws.Range("D1").Value = "politica"

# Row 2
$ws.Range("C2").Value = "Modelo identificado mas fora do range de preco"
$ws.Range("D2").Value = ""
$ws.Range("F2").Value = "classico"
$ws.Range("G2").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-bob-storm-90a-bivolt-automatico-cor-preto/p/MLB21562641?pdp_filters=seller_id:19365747#searchVariation=MLB21562641&position=1&search_layout=stack&type=product&tracking_id=b2671fdd-9475-4222-b116-ad49a565daa4"

# Row 3
$ws.Range("C3").Value = "FONTE 60A"
$ws.Range("D3").Value = "Acima"
$ws.Range("F3").Value = "premium"
$ws.Range("G3").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-60a-bivolt-storm-com-medidor-cca/p/MLB21320712?pdp_filters=seller_id:19365747#searchVariation=MLB21320712&position=2&search_layout=stack&type=product&tracking_id=b2671fdd-9475-4222-b116-ad49a565daa4"

# Row 4
$ws.Range("C4").Value = "FONTE 70A"
$ws.Range("D4").Value = "Igual"
$ws.Range("F4").Value = "classico"
$ws.Range("G4").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-70a-bivolt-com-medidor-cca/p/MLB21455208?pdp_filters=seller_id:19365747#searchVariation=MLB21455208&position=5&search_layout=stack&type=product&tracking_id=b2671fdd-9475-4222-b116-ad49a565daa4"

# Row 5
$ws.Range("C5").Value = "FONTE 40A"
$ws.Range("D5").Value = "Acima"
$ws.Range("F5").Value = "premium"
$ws.Range("G5").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-storm-40a-bivolt-12v-cor-preto/p/MLB22569833?pdp_filters=seller_id:19365747#searchVariation=MLB22569833&position=4&search_layout=stack&type=product&tracking_id=b2671fdd-9475-4222-b116-ad49a565daa4"

# Row 6
$ws.Range("C6").Value = "FONTE 120A"
$ws.Range("D6").Value = "Igual"
$ws.Range("F6").Value = "premium"
$ws.Range("G6").Value = "https://produto.mercadolivre.com.br/MLB-4423375714-fonte-digital-jfa-storm-automotiva-120a-bivolt-carregador-_JM#position%3D7%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Db2671fdd-9475-4222-b116-ad49a565daa4"

# Row 7
$ws.Range("C7").Value = "Sem Modelo"
$ws.Range("D7").Value = ""
$ws.Range("F7").Value = "classico"
$ws.Range("G7").Value = "https://www.mercadolivre.com.br/controle-longa-distncia-reposico-jfa-tx-k1200-som/p/MLB29541981?pdp_filters=seller_id:19365747#searchVariation=MLB29541981&position=6&search_layout=stack&type=product&tracking_id=b2671fdd-9475-4222-b116-ad49a565daa4"

# Row 8
$ws.Range("C8").Value = "Sem Modelo"
$ws.Range("D8").Value = ""
$ws.Range("F8").Value = "classico"
$ws.Range("G8").Value = "https://produto.mercadolivre.com.br/MLB-1948538513-controle-de-longa-distncia-som-automotivo-jfa-k1200-azul-_JM#position%3D8%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Db2671fdd-9475-4222-b116-ad49a565daa4"

# Row 9
$ws.Range("C9").Value = "FONTE 40A"
$ws.Range("D9").Value = "Igual"
$ws.Range("F9").Value = "classico"
$ws.Range("G9").Value = "https://produto.mercadolivre.com.br/MLB-4423372216-fonte-storm-digital-jfa-40a-carregador-com-medidor-de-cca-_JM#position%3D9%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Db2671fdd-9475-4222-b116-ad49a565daa4"

# Row 10
$ws.Range("C10").Value = "FONTE 70A"
$ws.Range("D10").Value = "Igual"
$ws.Range("F10").Value = "premium"
$ws.Range("G10").Value = "https://produto.mercadolivre.com.br/MLB-4423381104-carregador-jfa-digital-storm-fonte-automotiva-70a-medidor-_JM#position%3D10%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Db2671fdd-9475-4222-b116-ad49a565daa4"

# Row 11
$ws.Range("C11").Value = "FONTE 120A"
$ws.Range("D11").Value = "Igual"
$ws.Range("F11").Value = "classico"
$ws.Range("G11").Value = "https://produto.mercadolivre.com.br/MLB-4423388158-fonte-digital-jfa-storm-automotiva-carregador-120a-bivolt-_JM#position%3D11%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Db2671fdd-9475-4222-b116-ad49a565daa4"

# Row 12
$ws.Range("C12").Value = "FONTE 60A"
$ws.Range("D12").Value = "Acima"
$ws.Range("F12").Value = "classico"
$ws.Range("G12").Value = "https://produto.mercadolivre.com.br/MLB-3587094409-carregador-digital-fonte-jfa-storm-60a-medidor-de-cca-_JM#position%3D12%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Db2671fdd-9475-4222-b116-ad49a565daa4"

# Row 13
$ws.Range("C13").Value = "Sem Modelo"
$ws.Range("D13").Value = ""
$ws.Range("F13").Value = "classico"
$ws.Range("G13").Value = "https://produto.mercadolivre.com.br/MLB-1948540697-controle-de-longa-distncia-universal-jfa-k1200-azul-1200m-_JM#position%3D13%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Db2671fdd-9475-4222-b116-ad49a565daa4"

# Row 14
$ws.Range("C14").Value = "Sem Modelo"
$ws.Range("D14").Value = ""
$ws.Range("F14").Value = "premium"
$ws.Range("G14").Value = "https://produto.mercadolivre.com.br/MLB-3586968937-sequenciador-voltimetro-digital-jfa-altabaixa-voltagem-_JM#position%3D14%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Db2671fdd-9475-4222-b116-ad49a565daa4"

# Row 15
$ws.Range("C15").Value = "Sem Modelo"
$ws.Range("D15").Value = ""
$ws.Range("F15").Value = "premium"
$ws.Range("G15").Value = "https://produto.mercadolivre.com.br/MLB-3590198911-filtro-jfa-rca-anti-ruido-2020k-eletromagnetico-stereo-_JM#position%3D15%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Db2671fdd-9475-4222-b116-ad49a565daa4"

# Row 16
$ws.Range("C16").Value = "Sem Modelo"
$ws.Range("D16").Value = ""
$ws.Range("F16").Value = "classico"
$ws.Range("G16").Value = "https://produto.mercadolivre.com.br/MLB-3587032375-voltimetro-digital-jfa-sequenciador-altabaixa-voltagem-_JM#position%3D16%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Db2671fdd-9475-4222-b116-ad49a565daa4"

# Row 17
$ws.Range("C17").Value = "Sem Modelo"
$ws.Range("D17").Value = ""
$ws.Range("F17").Value = "classico"
$ws.Range("G17").Value = "https://produto.mercadolivre.com.br/MLB-4431471070-filtro-rca-anti-ruido-jfa-2020k-eletromagnetico-stereo-_JM#position%3D17%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Db2671fdd-9475-4222-b116-ad49a565daa4"

# Row 18
$ws.Range("C18").Value = "Sem Modelo"
$ws.Range("D18").Value = ""
$ws.Range("F18").Value = "premium"
$ws.Range("G18").Value = "https://produto.mercadolivre.com.br/MLB-1803852214-controle-longa-distncia-avulso-jfa-varias-cores-k1200-cx-_JM#position%3D18%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Db2671fdd-9475-4222-b116-ad49a565daa4"

# Row 19
$ws.Range("C19").Value = "Sem Modelo"
$ws.Range("D19").Value = ""
$ws.Range("F19").Value = "classico"
$ws.Range("G19").Value = "https://produto.mercadolivre.com.br/MLB-1803855509-controle-longa-distncia-avulso-jfa-k1200-cx-varias-cores-_JM#position%3D19%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Db2671fdd-9475-4222-b116-ad49a565daa4"

# Row 20
$ws.Range("C20").Value = "Sem Modelo"
$ws.Range("D20").Value = ""
$ws.Range("F20").Value = "premium"
$ws.Range("G20").Value = "https://produto.mercadolivre.com.br/MLB-2709379077-kit-4-controle-longa-distancia-jfa-azul-avulso-k1200-_JM#position%3D20%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Db2671fdd-9475-4222-b116-ad49a565daa4"

# Row 21
$ws.Range("C21").Value = "Sem Modelo"
$ws.Range("D21").Value = ""
$ws.Range("F21").Value = "premium"
$ws.Range("G21").Value = "https://produto.mercadolivre.com.br/MLB-2709347720-par-controle-longa-distancia-reposico-jfa-azul-k1200-avulso-_JM#position%3D21%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Db2671fdd-9475-4222-b116-ad49a565daa4"

# Row 22
$ws.Range("C22").Value = "Sem Modelo"
$ws.Range("D22").Value = ""
$ws.Range("F22").Value = "premium"
$ws.Range("G22").Value = "https://produto.mercadolivre.com.br/MLB-3587086063-conversor-rca-jfa-adaptador-comando-remoto-slim-automotivo-_JM#position%3D22%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Db2671fdd-9475-4222-b116-ad49a565daa4"

# Row 23
$ws.Range("C23").Value = "Sem Modelo"
$ws.Range("D23").Value = ""
$ws.Range("F23").Value = "classico"
$ws.Range("G23").Value = "https://produto.mercadolivre.com.br/MLB-2709355441-kit-com-4-controle-longa-distancia-jfa-azul-k1200-avulso-_JM#position%3D23%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Db2671fdd-9475-4222-b116-ad49a565daa4"

# Row 24
$ws.Range("C24").Value = "Sem Modelo"
$ws.Range("D24").Value = ""
$ws.Range("F24").Value = "classico"
$ws.Range("G24").Value = "https://produto.mercadolivre.com.br/MLB-2709334936-kit-2-controle-longa-distancia-azul-jfa-k1200-avulso-_JM#position%3D24%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Db2671fdd-9475-4222-b116-ad49a565daa4"

# Row 25
$ws.Range("C25").Value = "Sem Modelo"
$ws.Range("D25").Value = ""
$ws.Range("F25").Value = "premium"
$ws.Range("G25").Value = "https://produto.mercadolivre.com.br/MLB-3141326339-controle-reposico-longa-distncia-som-vermelho-jfa-tx-k1200-_JM#position%3D25%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Db2671fdd-9475-4222-b116-ad49a565daa4"
